{"js": "// \"wrong tag and remove cor. account\"\n// Remove the \"\u041a\u043e\u0440\u0440. \u0441\u0447\u0435\u0442: @<CORR_ACCOUNT>@\" paragraph (the correspondent-\n// account line in the bank requisites block) entirely. The paragraph right\n// above it (\"\u0411\u0418\u041a: @<BIK>@\") is left untouched.\nconst body = context.document.body;\n\nconst results = body.search(\"@<CORR_ACCOUNT>@\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\n// Delete every paragraph that contains the CORR_ACCOUNT placeholder (there is\n// exactly one in this document, but loop defensively in case of repeats).\nfor (let i = 0; i < results.items.length; i++) {\n  const paragraphs = results.items[i].paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < paragraphs.items.length; j++) {\n    paragraphs.items[j].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# \"wrong tag and remove cor. account\"\n# Remove the \"\u041a\u043e\u0440\u0440. \u0441\u0447\u0435\u0442: @<CORR_ACCOUNT>@\" paragraph (the correspondent-\n# account line in the bank requisites block) entirely. The paragraph right\n# above it (\"\u0411\u0418\u041a: @<BIK>@\") is left untouched.\n$d = $word.ActiveDocument\n\n# Collect every paragraph containing the CORR_ACCOUNT placeholder first (there\n# is exactly one in this document, but gather defensively in case of repeats),\n# then delete them back-to-front so earlier indices stay valid.\n$matches = @()\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs[$i]\n    if ($p.Range.Text -like \"*CORR_ACCOUNT*\") {\n        $matches += $p\n    }\n}\n\nfor ($k = $matches.Count - 1; $k -ge 0; $k--) {\n    $matches[$k].Range.Delete()\n}\n"}
